# D4H_Trackingsheet_v.1.0.xlsx update
#  - switch calculation to manual
#  - Doc sheet: move viewport/selection
#  - Product sheet: insert a "Start date" column, fill in new task #4
#    (Xac thuc tai khoan) rows, and make Product the active tab/sheet

$wb = $excel.ActiveWorkbook

# Workbook now opens in manual calculation mode
$excel.Calculation = -4135

$wsDoc = $wb.Worksheets.Item("Doc")
$wsProduct = $wb.Worksheets.Item("Product")

# ---- Product sheet: insert new "Start date" column before column E ----
$wsProduct.Columns("E:E").Insert()

$wsProduct.Range("E4").Value = "Start date"

$wsProduct.Range("E8").Value = 43758
$wsProduct.Range("E8").NumberFormat = "d-mmm"

# ---- Product sheet: new task #4 "Xac thuc tai khoan" (rows 13-16) ----
$wsProduct.Range("A13").Value = 4
$wsProduct.Range("B13").Value = "Xác thực tài khoản"
$wsProduct.Range("C13").Value = "Tạo chức năng nhập"
$wsProduct.Range("D13").Value = "Dương"
$wsProduct.Range("E13").Value = 43789
$wsProduct.Range("F13").Value = 43802
$wsProduct.Range("G13").Value = 43802
$wsProduct.Range("H13").Value = 43794
$wsProduct.Range("E13:H13").NumberFormat = "d-mmm"

$wsProduct.Range("C14").Value = "Tạo chức năng xuất"
$wsProduct.Range("D14").Value = "Nhật"
$wsProduct.Range("H14").Value = 43795
$wsProduct.Range("H14").NumberFormat = "d-mmm"

$wsProduct.Range("C15").Value = "Tạo chức năng đăng ký"
$wsProduct.Range("D15").Value = "Huyền"
$wsProduct.Range("H15").Value = 43796
$wsProduct.Range("H15").NumberFormat = "d-mmm"

$wsProduct.Range("C16").Value = "Test"
$wsProduct.Range("D16").Value = "Hoa, Hương"
$wsProduct.Range("H16").Value = 43797
$wsProduct.Range("H16").NumberFormat = "d-mmm"

# ---- Doc sheet viewport/selection ----
$wsDoc.Activate()
$wsDoc.Range("D9:D11").Select()

# ---- Product sheet becomes the active tab, with its own selection ----
$wsProduct.Activate()
$wsProduct.Range("G13").Select()
